$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("European design. & varieties")

# New rows 64-67: German wine regions added to the cepage names table.
# Values are written in the exact order the original authoring tool
# appended them to the shared-string table so new si entries line up.
$ws.Range("A64").Value = "Mosel-Saar-Ruwer (White), Germany"
$ws.Range("B64").Value = "Riesling"

$ws.Range("A65").Value = "Rheingau (White), Germany"
$ws.Range("B65").Value = "Riesling, Pinot noir"

$ws.Range("B66").Value = "Müller-Thurgau, Sylvaner"
$ws.Range("A66").Value = "Rheinhessen, Germany"

$ws.Range("A67").Value = "Pfalz, Germany"
$ws.Range("B67").Value = "Müller-Thurgau, Riesling, Sylvaner, Kerner, Scheurebe and Spätburgunder (Pinot Noir)"

# Move the view/selection to match the saved state in the workbook.
[void]$ws.Range("A12").Select()
